# Update gh-pages to output generated at 456a3b4
# Updates the "F" column (registration/view counts) on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 6969
$ws1.Range("F3").Value  = 100
$ws1.Range("F7").Value  = 6845
$ws1.Range("F13").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F20").Value = 5220
$ws1.Range("F21").Value = 0
$ws1.Range("F22").Value = 0
$ws1.Range("F23").Value = 635
$ws1.Range("F25").Value = 0

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 6969
$ws4.Range("F3").Value  = 0
$ws4.Range("F4").Value  = 0
$ws4.Range("F5").Value  = 0
$ws4.Range("F6").Value  = 0
$ws4.Range("F7").Value  = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F13").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F23").Value = 118
$ws4.Range("F25").Value = 635
$ws4.Range("F27").Value = 234
